# Generate Report for Handoff
# This script updates the localization-status workbook:
#  - Overview sheet: row order of the two files is effectively swapped, and the
#    status/date of the 6c124... file is updated to "Ready for handoff".
#  - zh-cn / de-de sheets: same row-content swap, status updated to
#    "Ready for handoff" for both rows, new handoff timestamps, and an error
#    detail message recorded for the 6c124... row.
#  - Column P (Error Detail) is widened on the zh-cn / de-de sheets to fit the
#    new long error message.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws1.Range("B2").Value = "e2e\d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws1.Range("A3").Value = "6c124c18-c336-4b27-8870-b61637c34677.md"
$ws1.Range("B3").Value = "e2e\6c124c18-c336-4b27-8870-b61637c34677.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-16 18:42:21"

# Hyperlinks on B2/B3 keep pointing at the same targets, but the displayed
# text is swapped to match the new A/B column content.
$hlTarget6c124 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/6c124c18-c336-4b27-8870-b61637c34677.md"
$hlTargetD156 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $hlTarget6c124, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\d156b7a4-8b41-4910-9674-238ac0bccee5.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), $hlTargetD156, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\6c124c18-c336-4b27-8870-b61637c34677.md")

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("G2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.d231fb23bcdc271149ae9cd7452341ec2fc6919e.zh-cn.xlf"
$ws2.Range("I2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws2.Range("J2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.d231fb23bcdc271149ae9cd7452341ec2fc6919e.zh-cn.xlf"
$ws2.Range("A3").Value = "6c124c18-c336-4b27-8870-b61637c34677.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "6c124c18-c336-4b27-8870-b61637c34677.f9a1e37d3206652964597145d91445a85125f3ce.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-16 18:42:17"
$ws2.Range("I3").Value = "6c124c18-c336-4b27-8870-b61637c34677.md"
$ws2.Range("J3").Value = "6c124c18-c336-4b27-8870-b61637c34677.f9a1e37d3206652964597145d91445a85125f3ce.zh-cn.xlf"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/6c124c18-c336-4b27-8870-b61637c34677.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0dd6de74cefbe77134140f6d8090ac631942b620/e2e/6c124c18-c336-4b27-8870-b61637c34677.md."

$ws2HlTarget6c124 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/6c124c18-c336-4b27-8870-b61637c34677.md"
$ws2HlTarget6c124Zhcn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/77a08c849a2dfba43ea456aec716f50e986927e5/e2e/6c124c18-c336-4b27-8870-b61637c34677.md"
$ws2HlTargetD156 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws2HlTargetD156Zhcn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/77a08c849a2dfba43ea456aec716f50e986927e5/e2e/d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2HlTarget6c124, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "d156b7a4-8b41-4910-9674-238ac0bccee5.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $ws2HlTarget6c124Zhcn, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "d156b7a4-8b41-4910-9674-238ac0bccee5.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2HlTargetD156, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "6c124c18-c336-4b27-8870-b61637c34677.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $ws2HlTargetD156Zhcn, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "6c124c18-c336-4b27-8870-b61637c34677.md")

# Widen the Error Detail column to fit the new long message.
$ws2.Columns.Item(16).ColumnWidth = 39.1667

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("G2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.d231fb23bcdc271149ae9cd7452341ec2fc6919e.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-16 18:42:21"
$ws3.Range("I2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws3.Range("J2").Value = "d156b7a4-8b41-4910-9674-238ac0bccee5.d231fb23bcdc271149ae9cd7452341ec2fc6919e.de-de.xlf"
$ws3.Range("A3").Value = "6c124c18-c336-4b27-8870-b61637c34677.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "6c124c18-c336-4b27-8870-b61637c34677.f9a1e37d3206652964597145d91445a85125f3ce.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-16 18:42:21"
$ws3.Range("I3").Value = "6c124c18-c336-4b27-8870-b61637c34677.md"
$ws3.Range("J3").Value = "6c124c18-c336-4b27-8870-b61637c34677.f9a1e37d3206652964597145d91445a85125f3ce.de-de.xlf"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/6c124c18-c336-4b27-8870-b61637c34677.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0dd6de74cefbe77134140f6d8090ac631942b620/e2e/6c124c18-c336-4b27-8870-b61637c34677.md."

$ws3HlTarget6c124 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/6c124c18-c336-4b27-8870-b61637c34677.md"
$ws3HlTarget6c124Dede = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/df4370cf21cdec76d4cb59dc860ac60300de4c76/e2e/6c124c18-c336-4b27-8870-b61637c34677.md"
$ws3HlTargetD156 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c07c96814a2c2f23b51538829958056fddc6f8fc/e2e/d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws3HlTargetD156Dede = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/df4370cf21cdec76d4cb59dc860ac60300de4c76/e2e/d156b7a4-8b41-4910-9674-238ac0bccee5.md"
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3HlTarget6c124, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "d156b7a4-8b41-4910-9674-238ac0bccee5.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $ws3HlTarget6c124Dede, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "d156b7a4-8b41-4910-9674-238ac0bccee5.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ws3HlTargetD156, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "6c124c18-c336-4b27-8870-b61637c34677.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $ws3HlTargetD156Dede, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "6c124c18-c336-4b27-8870-b61637c34677.md")

# Widen the Error Detail column to fit the new long message.
$ws3.Columns.Item(16).ColumnWidth = 39.1667
